# Auto-generated edit script: update column F (报名人数/浏览量等) values per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 996  # F2: 995 -> 996
$ws.Cells.Item(4, 6).Value = 1229  # F4: 1227 -> 1229
$ws.Cells.Item(5, 6).Value = 49  # F5: 47 -> 49
$ws.Cells.Item(6, 6).Value = 728  # F6: 727 -> 728
$ws.Cells.Item(7, 6).Value = 1069  # F7: 1068 -> 1069
$ws.Cells.Item(8, 6).Value = 4575  # F8: 4567 -> 4575
$ws.Cells.Item(9, 6).Value = 591  # F9: 588 -> 591
$ws.Cells.Item(10, 6).Value = 160  # F10: 159 -> 160
$ws.Cells.Item(11, 6).Value = 1750  # F11: 1749 -> 1750
$ws.Cells.Item(12, 6).Value = 29  # F12: 27 -> 29
$ws.Cells.Item(13, 6).Value = 691  # F13: 682 -> 691
$ws.Cells.Item(14, 6).Value = 31  # F14: 30 -> 31
$ws.Cells.Item(17, 6).Value = 1109  # F17: 1106 -> 1109
$ws.Cells.Item(18, 6).Value = 1551  # F18: 1550 -> 1551
$ws.Cells.Item(19, 6).Value = 791  # F19: 786 -> 791
$ws.Cells.Item(20, 6).Value = 701  # F20: 696 -> 701
$ws.Cells.Item(21, 6).Value = 531  # F21: 530 -> 531
$ws.Cells.Item(24, 6).Value = 113  # F24: 109 -> 113
$ws.Cells.Item(27, 6).Value = 364  # F27: 361 -> 364
$ws.Cells.Item(28, 6).Value = 2484  # F28: 2483 -> 2484
$ws.Cells.Item(29, 6).Value = 288  # F29: 287 -> 288
$ws.Cells.Item(30, 6).Value = 1492  # F30: 1483 -> 1492
$ws.Cells.Item(32, 6).Value = 14  # F32: 13 -> 14
$ws.Cells.Item(34, 6).Value = 4143  # F34: 4135 -> 4143

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 9  # F4: 10 -> 9
$ws.Cells.Item(5, 6).Value = 210  # F5: 209 -> 210
$ws.Cells.Item(7, 6).Value = 15  # F7: 14 -> 15
$ws.Cells.Item(12, 6).Value = 355  # F12: 351 -> 355
$ws.Cells.Item(13, 6).Value = 4152  # F13: 4150 -> 4152
$ws.Cells.Item(18, 6).Value = 12  # F18: 1 -> 12
$ws.Cells.Item(21, 6).Value = 273  # F21: 271 -> 273
$ws.Cells.Item(25, 6).Value = 136  # F25: 134 -> 136
$ws.Cells.Item(31, 6).Value = 1736  # F31: 1735 -> 1736

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 1314  # F4: 1312 -> 1314
$ws.Cells.Item(5, 6).Value = 1713  # F5: 1711 -> 1713
$ws.Cells.Item(7, 6).Value = 1061  # F7: 1058 -> 1061
$ws.Cells.Item(8, 6).Value = 180  # F8: 173 -> 180

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1314  # F2: 1312 -> 1314
$ws.Cells.Item(3, 6).Value = 1713  # F3: 1711 -> 1713
$ws.Cells.Item(5, 6).Value = 1061  # F5: 1058 -> 1061
$ws.Cells.Item(6, 6).Value = 996  # F6: 995 -> 996
$ws.Cells.Item(7, 6).Value = 1229  # F7: 1227 -> 1229
$ws.Cells.Item(9, 6).Value = 49  # F9: 47 -> 49
$ws.Cells.Item(10, 6).Value = 728  # F10: 727 -> 728
$ws.Cells.Item(11, 6).Value = 210  # F11: 209 -> 210
$ws.Cells.Item(12, 6).Value = 210  # F12: 209 -> 210
$ws.Cells.Item(13, 6).Value = 180  # F13: 173 -> 180
$ws.Cells.Item(14, 6).Value = 15  # F14: 14 -> 15
$ws.Cells.Item(15, 6).Value = 1069  # F15: 1068 -> 1069
$ws.Cells.Item(17, 6).Value = 4575  # F17: 4567 -> 4575
$ws.Cells.Item(18, 6).Value = 591  # F18: 588 -> 591
$ws.Cells.Item(19, 6).Value = 160  # F19: 159 -> 160
$ws.Cells.Item(20, 6).Value = 1750  # F20: 1749 -> 1750
$ws.Cells.Item(21, 6).Value = 29  # F21: 27 -> 29
$ws.Cells.Item(22, 6).Value = 691  # F22: 682 -> 691
$ws.Cells.Item(23, 6).Value = 355  # F23: 351 -> 355
$ws.Cells.Item(27, 6).Value = 1109  # F27: 1106 -> 1109
$ws.Cells.Item(28, 6).Value = 1551  # F28: 1550 -> 1551
$ws.Cells.Item(31, 6).Value = 791  # F31: 786 -> 791
$ws.Cells.Item(32, 6).Value = 701  # F32: 696 -> 701
$ws.Cells.Item(33, 6).Value = 531  # F33: 530 -> 531
$ws.Cells.Item(36, 6).Value = 113  # F36: 109 -> 113
$ws.Cells.Item(38, 6).Value = 273  # F38: 271 -> 273
$ws.Cells.Item(43, 6).Value = 364  # F43: 361 -> 364
$ws.Cells.Item(44, 6).Value = 2484  # F44: 2483 -> 2484
$ws.Cells.Item(46, 6).Value = 1492  # F46: 1483 -> 1492
$ws.Cells.Item(48, 6).Value = 14  # F48: 13 -> 14
$ws.Cells.Item(49, 6).Value = 4143  # F49: 4135 -> 4143
